# Movie_Data_Template.xlsx edit:
#  - Remove "International Gross" column (column E) from MovieInfo sheet
#  - Rename "MovieCrew" sheet to "Director's Highest gross Films"
#  - Update view/selection state to match final saved state

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("MovieInfo")
$ws2 = $wb.Worksheets.Item("MovieCrew")

# Remove the "International Gross" column (column E) on MovieInfo.
# This shifts F:I left to E:H and automatically fixes up the
# shared-string table and any references to it (e.g. on MovieCrew).
$ws1.Columns("E").Delete() | Out-Null

# Rename the crew sheet.
$ws2.Name = "Director's Highest gross Films"

# Restore per-sheet selection/view state seen in the final workbook:
# MovieCrew (now renamed) keeps a stored selection at D26 ...
$ws2.Activate() | Out-Null
$ws2.Range("D26").Select() | Out-Null

# ... while MovieInfo ends up as the active/selected sheet with
# column E selected end-to-end.
$ws1.Activate() | Out-Null
$ws1.Range("E1:E1048576").Select() | Out-Null
